$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: C1 keeps "CompanyID"; D1 becomes "Headline"; add E1 "Summary", F1 "Language"
$ws.Range("D1").Value = "Headline"
$ws.Range("E1").Value = "Summary"
$ws.Range("F1").Value = "Language"

# Make C1:F1 bold header styling (match existing bold header cells)
$ws.Range("C1:F1").Font.Bold = $true

# Update data row values
$ws.Range("A2").Value = 45595
$ws.Range("C2").Value = "DEMC"

# Update selection to match the new active cell
$ws.Range("F2").Select()
